# Lipsync-Commands workbook update:
#  - Wireless sheet: reorder the Bluetooth rows, change "BT,1:{N}" -> "BT,1:3",
#    insert two new "Communication mode" rows, and move "Perform factory reset"
#    to the end of the table.
#  - Make the Wireless sheet the active/selected tab (it was Mouse before).

$wb  = $excel.ActiveWorkbook
$wsWireless = $wb.Worksheets.Item("Wireless")

# --- Wireless sheet: restructure rows 41-43 -------------------------------
# Before:
#   41: FR,0:0  | SUCCESS:FR,0:0                       | FAIL:SETTINGS | Perform factory reset
#   42: BT,0:0  | SUCCESS:BT,0:{Bluetooth Module Mode}  | FAIL:SETTINGS | Get Bluetooth module mode (...)
#   43: BT,1:{N}| SUCCESS:BT,1:{Bluetooth Module Mode}  | FAIL:SETTINGS | Set Bluetooth module mode (...)
#
# After:
#   41: BT,0:0  | SUCCESS:BT,0:{Bluetooth Module Mode}  | FAIL:SETTINGS | Get Bluetooth module mode (...)
#   42: BT,1:3  | SUCCESS:BT,1:{Bluetooth Module Mode}  | FAIL:SETTINGS | Set Bluetooth module mode (...)
#   43: CM,0:0  | SUCCESS:CM,0:{N}                      | FAIL:SETTINGS | Get Communication mode (0=USB , Bluetooth =1)
#   44: CM,1:{N}| SUCCESS:CM,1:{N}                      | FAIL:SETTINGS | Set Communication mode (0=USB , Bluetooth =1)
#   45: FR,0:0  | SUCCESS:FR,0:0                        | FAIL:SETTINGS | Perform factory reset

# 1) Remove the old "Factory reset" row (41) -- this shifts the two BT rows
#    up to 41/42.
$wsWireless.Rows.Item(41).Delete()

# 2) The old "BT,1:{N}" command (now at row 42) becomes "BT,1:3".
$wsWireless.Range("A42").Value = "BT,1:3"

# 3) Insert two fresh rows for the new "Communication mode" commands right
#    after the Bluetooth rows (at 43 and 44).
$wsWireless.Rows.Item(43).Insert()
$wsWireless.Rows.Item(43).Insert()

$wsWireless.Range("A43").Value = "CM,0:0"
$wsWireless.Range("B43").Value = "SUCCESS:CM,0:{N}"
$wsWireless.Range("C43").Value = "FAIL:SETTINGS"
$wsWireless.Range("D43").Value = "Get Communication mode (0=USB , Bluetooth =1)"

$wsWireless.Range("A44").Value = "CM,1:{N}"
$wsWireless.Range("B44").Value = "SUCCESS:CM,1:{N}"
$wsWireless.Range("C44").Value = "FAIL:SETTINGS"
$wsWireless.Range("D44").Value = "Set Communication mode (0=USB , Bluetooth =1)"

# 4) Re-append the "Perform factory reset" row at the end of the table (45).
$wsWireless.Range("A45").Value = "FR,0:0"
$wsWireless.Range("B45").Value = "SUCCESS:FR,0:0"
$wsWireless.Range("C45").Value = "FAIL:SETTINGS"
$wsWireless.Range("D45").Value = "Perform factory reset"

# --- Selection / active tab -----------------------------------------------
# The Wireless sheet becomes the active tab, with D44 selected.
$wsWireless.Range("D44").Select() | Out-Null
$wsWireless.Activate()
